$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.342.11"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "3.503.98"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.56"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.50"
$ws.Range("E6").Value = "  +2.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.487"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.388"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "4.094.83"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000183"
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("D15").Value = "3.497.17"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.08"
$ws.Range("E16").Value = "  -4.54%  "
$ws.Range("D17").Value = "64.315.43"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.93"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.76"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.76"
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "395.14"
$ws.Range("E21").Value = "  +3.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.574"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "3.639.33"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.26"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.66"
$ws.Range("E27").Value = "  +2.18%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.44"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -5.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.34"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "3.518.01"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  +4.35%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.51"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.20"
$ws.Range("E37").Value = "  -2.46%  "
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.93"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "161.41"
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0785"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.808"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.17"
$ws.Range("E44").Value = "  -4.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.44"
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.66"
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.18"
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("D48").Value = "2.469.48"
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.80"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.898"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("E51").Value = "  -1.11%  "
